$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "27.030.44"
$ws.Range("D3").Value = "1.675.22"
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "215.30"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.25%  "
Set-TextValue "D9" "21.43"
$ws.Range("E9").Value = "  +5.66%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "1.913.05"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "1.693.19"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("E14").Value = "  +0.83%  "
Set-TextValue "D15" "0.535"
$ws.Range("E15").Value = "  +1.69%  "
Set-TextValue "D16" "66.34"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "27.036.71"
$ws.Range("E17").Value = "  +0.58%  "
Set-TextValue "D18" "8.17"
$ws.Range("E18").Value = "  +2.74%  "
Set-TextValue "D19" "236.10"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +1.67%  "
Set-TextValue "D23" "9.26"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  -1.94%  "
Set-TextValue "D25" "148.17"
$ws.Range("E25").Value = "  +1.03%  "
Set-TextValue "D26" "7.27"
$ws.Range("E26").Value = "  +2.19%  "
Set-TextValue "D27" "16.50"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  +0.01%  "
Set-TextValue "D30" "0.0498"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").Value = "1.542.40"
$ws.Range("E33").Value = "  +6.64%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E35").Value = "  +5.07%  "
Set-TextValue "D36" "2.38"
$ws.Range("E36").Value = "  -1.13%  "
Set-TextValue "D37" "0.591"
$ws.Range("E37").Value = "  +1.07%  "
Set-TextValue "D38" "0.912"
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("E39").Value = "  +2.05%  "
Set-TextValue "D40" "1.04"
$ws.Range("E40").Value = "  +4.19%  "
$ws.Range("E41").Value = "  +0.05%  "
Set-TextValue "D42" "67.78"
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "1.820.35"
$ws.Range("E45").Value = "  +0.96%  "
Set-TextValue "D46" "0.779"
$ws.Range("E46").Value = "  -0.25%  "
Set-TextValue "D47" "90.45"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  +2.18%  "
Set-TextValue "D50" "8.05"
$ws.Range("E50").Value = "  +6.35%  "
$ws.Range("E51").Value = "  +0.31%  "

Write-Output "Updated 74 cells"
